$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 55
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 45
